# Auto-generated script to apply 2024-08-14 data updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 4954
$ws.Range('K3').Value = 5091
$ws.Range('J4').Value = 1831
$ws.Range('K4').Value = 1055
$ws.Range('K5').Value = 360
$ws.Range('K6').Value = 5712
$ws.Range('J7').Value = 29299
$ws.Range('K7').Value = 17172

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 147
$ws.Range('K4').Value = 63
$ws.Range('K6').Value = 129
$ws.Range('K7').Value = 511
$ws.Range('K8').Value = 1149
$ws.Range('K11').Value = 334
$ws.Range('K14').Value = 92
$ws.Range('K15').Value = 171
$ws.Range('J19').Value = 858
$ws.Range('K19').Value = 519
$ws.Range('K20').Value = 393
$ws.Range('K23').Value = 173
$ws.Range('K24').Value = 52
$ws.Range('K25').Value = 82
$ws.Range('K29').Value = 921
$ws.Range('K33').Value = 726
$ws.Range('K34').Value = 91
$ws.Range('K36').Value = 224
$ws.Range('K37').Value = 581
$ws.Range('K42').Value = 637
$ws.Range('J44').Value = 231
$ws.Range('K48').Value = 215
$ws.Range('K49').Value = 96
$ws.Range('K51').Value = 216
$ws.Range('K52').Value = 447
$ws.Range('K53').Value = 223
$ws.Range('K54').Value = 337
$ws.Range('K63').Value = 49
$ws.Range('K65').Value = 390
$ws.Range('K66').Value = 56
$ws.Range('K67').Value = 663
$ws.Range('K70').Value = 32
$ws.Range('K77').Value = 123
$ws.Range('K79').Value = 419
$ws.Range('K81').Value = 11
$ws.Range('K83').Value = 377
$ws.Range('K84').Value = 130
$ws.Range('K85').Value = 789
$ws.Range('K86').Value = 115
$ws.Range('K89').Value = 249
$ws.Range('K91').Value = 185
$ws.Range('K92').Value = 65
$ws.Range('K93').Value = 65
$ws.Range('K94').Value = 227
$ws.Range('K96').Value = 183
$ws.Range('K97').Value = 136
$ws.Range('K99').Value = 294
$ws.Range('J101').Value = 29299
$ws.Range('K101').Value = 17172

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('K2').Value = 33
$ws.Range('K7').Value = 92

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K3').Value = 34
$ws.Range('K7').Value = 183

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K6').Value = 132
$ws.Range('K7').Value = 511

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K3').Value = 88
$ws.Range('K7').Value = 334

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K2').Value = 69
$ws.Range('K7').Value = 249

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 266
$ws.Range('K3').Value = 266
$ws.Range('K4').Value = 47
$ws.Range('K6').Value = 187
$ws.Range('K7').Value = 789

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K6').Value = 165
$ws.Range('K7').Value = 447

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K2').Value = 55
$ws.Range('K3').Value = 54
$ws.Range('K4').Value = 10
$ws.Range('K7').Value = 223

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 319
$ws.Range('K3').Value = 341
$ws.Range('K5').Value = 33
$ws.Range('K6').Value = 391
$ws.Range('K7').Value = 1149

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K3').Value = 138
$ws.Range('K6').Value = 87
$ws.Range('K7').Value = 377

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 199
$ws.Range('K3').Value = 270
$ws.Range('K6').Value = 209
$ws.Range('K7').Value = 726

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K3').Value = 190
$ws.Range('K6').Value = 175
$ws.Range('K7').Value = 581

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K2').Value = 118
$ws.Range('K6').Value = 151
$ws.Range('K7').Value = 390

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K6').Value = 74
$ws.Range('K7').Value = 294

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 191
$ws.Range('K3').Value = 231
$ws.Range('K6').Value = 189
$ws.Range('K7').Value = 663

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K3').Value = 51
$ws.Range('K7').Value = 130

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K4').Value = 10
$ws.Range('K7').Value = 96

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K2').Value = 52
$ws.Range('K6').Value = 178
$ws.Range('K7').Value = 337

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 267
$ws.Range('K3').Value = 329
$ws.Range('K6').Value = 254
$ws.Range('K7').Value = 921

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K2').Value = 29
$ws.Range('K7').Value = 215

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J4').Value = 44
$ws.Range('K6').Value = 163
$ws.Range('J7').Value = 858
$ws.Range('K7').Value = 519

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('J4').Value = 13
$ws.Range('J7').Value = 231

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('K3').Value = 35
$ws.Range('K6').Value = 37
$ws.Range('K7').Value = 129

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K4').Value = 24
$ws.Range('K6').Value = 243
$ws.Range('K7').Value = 637

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('K4').Value = 6
$ws.Range('K7').Value = 52

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K4').Value = 11
$ws.Range('K7').Value = 173

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K6').Value = 44
$ws.Range('K7').Value = 185

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K2').Value = 139
$ws.Range('K7').Value = 419

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K3').Value = 126
$ws.Range('K7').Value = 393

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K2').Value = 89
$ws.Range('K6').Value = 51
$ws.Range('K7').Value = 224

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('K6').Value = 25
$ws.Range('K7').Value = 65

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K2').Value = 32
$ws.Range('K7').Value = 91

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K6').Value = 96
$ws.Range('K7').Value = 227

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('K4').Value = 6
$ws.Range('K7').Value = 82

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K6').Value = 52
$ws.Range('K7').Value = 171

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('K2').Value = 16
$ws.Range('K7').Value = 56

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K2').Value = 48
$ws.Range('K7').Value = 147

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K6').Value = 81
$ws.Range('K7').Value = 136

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('K6').Value = 33
$ws.Range('K7').Value = 65

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range('K3').Value = 7
$ws.Range('K7').Value = 32

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('K3').Value = 20
$ws.Range('K7').Value = 115

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K2').Value = 62
$ws.Range('K7').Value = 216

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K3').Value = 50
$ws.Range('K7').Value = 123

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('K4').Value = 5
$ws.Range('K7').Value = 63

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range('K6').Value = 5
$ws.Range('K7').Value = 11
